$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column O (rows 5-8) with "FALSE" text values, matching the style of
# the plain body cells (column A). We paste the *value* from an existing
# "FALSE" text cell (column K, which already stores FALSE as literal text)
# and then paste the *format* from a plain body cell (column A) on top, so
# the new cells end up with exactly the same style/type as the rest of the
# table (no stray intermediate styles).
$ws.Range("K5").Copy()
$ws.Range("O5").PasteSpecial(-4163)
$ws.Range("A5").Copy()
$ws.Range("O5").PasteSpecial(-4122)

$ws.Range("K6").Copy()
$ws.Range("O6").PasteSpecial(-4163)
$ws.Range("A5").Copy()
$ws.Range("O6").PasteSpecial(-4122)

$ws.Range("K7").Copy()
$ws.Range("O7").PasteSpecial(-4163)
$ws.Range("A5").Copy()
$ws.Range("O7").PasteSpecial(-4122)

$ws.Range("K8").Copy()
$ws.Range("O8").PasteSpecial(-4163)
$ws.Range("A5").Copy()
$ws.Range("O8").PasteSpecial(-4122)

# Add new header cell "Internal Assignment" in O4, styled bold like the
# other headers (K4/L4/M4/N4) but with a larger, size-12 font.
$ws.Range("O4").Value = "Internal Assignment"
$ws.Range("O4").Font.Bold = $true
$ws.Range("O4").Font.Size = 12
$ws.Range("O4").Font.Name = "Calibri"
$ws.Range("O4").Font.Color = 0

# Update the active selection to match the target workbook state.
$ws.Range("O7:O8").Select()
